# Update cryptocurrency price/volume data in the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.964.18"
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").Value = "3.171.72"
$ws.Range("E3").Value = "  -4.13%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'591.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.65%  "
$ws.Range("D6").Value = "'134.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.01%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.169.28"
$ws.Range("E8").Value = "  -4.16%  "
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("E10").Value = "  -6.47%  "
$ws.Range("D11").Value = "'5.21"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.15%  "
$ws.Range("D12").Value = "'0.452"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.79%  "
$ws.Range("E13").Value = "  -4.78%  "
$ws.Range("D14").Value = "'34.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.05%  "
$ws.Range("D15").Value = "3.690.95"
$ws.Range("E15").Value = "  -4.23%  "
$ws.Range("E16").Value = "  -1.62%  "
$ws.Range("D17").Value = "3.169.71"
$ws.Range("E17").Value = "  -4.15%  "
$ws.Range("D18").Value = "62.937.64"
$ws.Range("E18").Value = "  -1.44%  "
$ws.Range("E19").Value = "  -5.13%  "
$ws.Range("D20").Value = "'459.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.47%  "
$ws.Range("D21").Value = "'13.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.47%  "
$ws.Range("E23").Value = "  -5.37%  "
$ws.Range("D24").Value = "'13.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.59%  "
$ws.Range("D25").Value = "'82.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.63%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("E28").Value = "  -4.27%  "
$ws.Range("D29").Value = "'6.72"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.94%  "
$ws.Range("D30").Value = "'7.62"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.92%  "
$ws.Range("E31").Value = "  -5.91%  "
$ws.Range("D32").Value = "'27.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.99%  "
$ws.Range("E33").Value = "  -4.75%  "
$ws.Range("D34").Value = "'2.36"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.61%  "
$ws.Range("D35").Value = "'1.02"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.68%  "
$ws.Range("E36").Value = "  -4.70%  "
$ws.Range("D37").Value = "'51.12"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.37%  "
$ws.Range("D38").Value = "0.0₃0705"
$ws.Range("E38").Value = "  -6.01%  "
$ws.Range("E39").Value = "  -3.43%  "
$ws.Range("D40").Value = "'402.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.94%  "
$ws.Range("E41").Value = "  -3.11%  "
$ws.Range("E43").Value = "  -6.49%  "
$ws.Range("D44").Value = "2.813.67"
$ws.Range("E44").Value = "  -9.70%  "
$ws.Range("E45").Value = "  -5.57%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("E47").Value = "  -6.02%  "
$ws.Range("D48").Value = "'123.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.69%  "
$ws.Range("D49").Value = "'34.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.06%  "
$ws.Range("D50").Value = "'25.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.15%  "
$ws.Range("E51").Value = "  -2.56%  "
